$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 5, shifting existing rows 5-15 down to 6-16.
$ws.Rows("5:5").Insert()

# Populate the newly inserted row 5 with the new data point.
# Column A keeps the date number format used by the rest of the column (style index 1).
$ws.Range("A5").Value = 42810
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 1

# Update the active selection to match the target workbook state.
$ws.Range("C6").Select()
